$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1, J1 - match style/formatting of existing header cells (e.g. H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data values for columns I and J, rows 2-18
$values = @{
    2  = @(7, 7)
    3  = @(7, 7)
    4  = @(5, 5)
    5  = @(5, 6)
    6  = @(9, 9)
    7  = @(4, 5)
    8  = @(7, 7)
    9  = @(2, 4)
    10 = @(8, 9)
    11 = @(7, 8)
    12 = @(5, 6)
    13 = @(6, 7)
    14 = @(6, 7)
    15 = @(7, 8)
    16 = @(1, 4)
    17 = @(10, 10)
    18 = @(3, 3)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
